# Update the "Translatable_Site_labels" sheet to add the new "edit page for texts"
# related labels (Admin section, submit/edit/undo buttons, country-description
# strings) and drop two labels that are no longer used.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translatable_Site_labels")

# --- 1. Remove the old duplicate "aka" row for the Author component (row 2) ---
$ws.Rows.Item(2).Delete()

# --- 2. Insert 5 new rows after the Author block (now ending at row 13,
#        "unspecified : 'not specified'") for the new author-edit-page labels ---
$ws.Range("A14:A18").EntireRow.Insert()

$ws.Range("B14").Value = "submit_edits"
$ws.Range("C14").Value = "Submit Edits"
$ws.Range("D14").Value = "Author"
$ws.Range("E14").Formula = '=_xlfn.CONCAT("",B14," : ''",C14,"'',")'

$ws.Range("B15").Value = "edit_country_birth_description"
$ws.Range("C15").Value = " (city, country) of birth"
$ws.Range("D15").Value = "Author"
$ws.Range("E15").Formula = '=_xlfn.CONCAT("",B15," : ''",C15,"'',")'

$ws.Range("B16").Value = "edit_country_death_description"
$ws.Range("C16").Value = " (city, country) of death"
$ws.Range("D16").Value = "Author"
$ws.Range("E16").Formula = '=_xlfn.CONCAT("",B16," : ''",C16,"'',")'

$ws.Range("B17").Value = "editBtn"
$ws.Range("C17").Value = "Edit"
$ws.Range("D17").Value = "Author"
$ws.Range("E17").Formula = '=_xlfn.CONCAT("",B17," : ''",C17,"'',")'

$ws.Range("B18").Value = "exitEditBtn"
$ws.Range("C18").Value = "Stop Editing"
$ws.Range("D18").Value = "Author"
$ws.Range("E18").Formula = '=_xlfn.CONCAT("",B18," : ''",C18,"'',")'

# --- 3. Remove the duplicate "unspecified : 'not specified'" row for the Text
#        component (now at row 28) ---
$ws.Rows.Item(28).Delete()

# --- 4. Remove the duplicate "isbn : 'ISBN10/13'" row for the Edition
#        component (now at row 36) ---
$ws.Rows.Item(36).Delete()

# --- 5. Append the new Import / Admin / Author-Text-Editions rows at the
#        bottom of the table (rows 48-56) ---
$ws.Range("B48").Value = "import_refresh"
$ws.Range("C48").Value = "Refresh"
$ws.Range("D48").Value = "Import"
$ws.Range("E48").Formula = '=_xlfn.CONCAT("",B48," : ''",C48,"'',")'

$ws.Range("B49").Value = "import_type_authors"
$ws.Range("C49").Value = "Authors"
$ws.Range("D49").Value = "Admin"
$ws.Range("E49").Formula = '=_xlfn.CONCAT("",B49," : ''",C49,"'',")'

$ws.Range("B50").Value = "import_type_texts"
$ws.Range("C50").Value = "Texts"
$ws.Range("D50").Value = "Admin"
$ws.Range("E50").Formula = '=_xlfn.CONCAT("",B50," : ''",C50,"'',")'

$ws.Range("B51").Value = "import_type_editions"
$ws.Range("C51").Value = "Editions"
$ws.Range("D51").Value = "Admin"
$ws.Range("E51").Formula = '=_xlfn.CONCAT("",B51," : ''",C51,"'',")'

$ws.Range("B52").Value = "import_error"
$ws.Range("C52").Value = "Data has not been imported or the data imported is empty"
$ws.Range("D52").Value = "Admin"
$ws.Range("E52").Formula = '=_xlfn.CONCAT("",B52," : ''",C52,"'',")'

$ws.Range("B53").Value = "import_databtn"
$ws.Range("C53").Value = "Imported data"
$ws.Range("D53").Value = "Admin"
$ws.Range("E53").Formula = '=_xlfn.CONCAT("",B53," : ''",C53,"'',")'

$ws.Range("B54").Value = "latest_editsbtn"
$ws.Range("C54").Value = "Latest edits"
$ws.Range("D54").Value = "Admin"
$ws.Range("E54").Formula = '=_xlfn.CONCAT("",B54," : ''",C54,"'',")'

$ws.Range("B55").Value = "admin_header"
$ws.Range("C55").Value = "Admin"
$ws.Range("D55").Value = "Admin"
$ws.Range("E55").Formula = '=_xlfn.CONCAT("",B55," : ''",C55,"'',")'

$ws.Range("B56").Value = "undoEditBtn"
$ws.Range("C56").Value = "Undo all changes"
$ws.Range("D56").Value = "Author/Text/Editions"
$ws.Range("E56").Formula = '=_xlfn.CONCAT("",B56," : ''",C56,"'',")'

# --- 6. Re-purpose the header cell E1 and add the closing JS snippet lines
#        E57 / E59 (row 58 intentionally left blank) so column E reads as a
#        complete `labels.js` export ---
$ws.Range("E1").Value = "const labels = {"
$ws.Range("E57").Value = "}"
$ws.Range("E59").Value = "export default labels;"

# --- 7. Restore the view: scroll down to the newly added rows and select E40 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E40").Select()
